$wb = $excel.ActiveWorkbook

# ---- PIR sheet: add rows 160-172 ----
$ws = $wb.Worksheets.Item("PIR")
$dateRng = $ws.Range("A160:A172")
$dateRng.NumberFormat = "@"
$ws.Cells.Item(160, 1).Value = "2026-01-28"
$ws.Cells.Item(160, 2).Value = "16:23:45"
$ws.Cells.Item(160, 3).Value = "16:00"
$ws.Cells.Item(160, 4).Value = "Bathroom"
$ws.Cells.Item(160, 5).Value = "No Motion"
$ws.Cells.Item(160, 6).Value = "Inactive"
$ws.Cells.Item(161, 1).Value = "2026-01-28"
$ws.Cells.Item(161, 2).Value = "16:23:47"
$ws.Cells.Item(161, 3).Value = "16:00"
$ws.Cells.Item(161, 4).Value = "Bathroom"
$ws.Cells.Item(161, 5).Value = "No Motion"
$ws.Cells.Item(161, 6).Value = "Inactive"
$ws.Cells.Item(162, 1).Value = "2026-01-28"
$ws.Cells.Item(162, 2).Value = "16:23:52"
$ws.Cells.Item(162, 3).Value = "16:00"
$ws.Cells.Item(162, 4).Value = "Bathroom"
$ws.Cells.Item(162, 5).Value = "No Motion"
$ws.Cells.Item(162, 6).Value = "Inactive"
$ws.Cells.Item(163, 1).Value = "2026-01-28"
$ws.Cells.Item(163, 2).Value = "16:23:57"
$ws.Cells.Item(163, 3).Value = "16:00"
$ws.Cells.Item(163, 4).Value = "Bathroom"
$ws.Cells.Item(163, 5).Value = "No Motion"
$ws.Cells.Item(163, 6).Value = "Inactive"
$ws.Cells.Item(164, 1).Value = "2026-01-28"
$ws.Cells.Item(164, 2).Value = "16:24:02"
$ws.Cells.Item(164, 3).Value = "16:00"
$ws.Cells.Item(164, 4).Value = "Bathroom"
$ws.Cells.Item(164, 5).Value = "No Motion"
$ws.Cells.Item(164, 6).Value = "Inactive"
$ws.Cells.Item(165, 1).Value = "2026-01-28"
$ws.Cells.Item(165, 2).Value = "16:24:07"
$ws.Cells.Item(165, 3).Value = "16:00"
$ws.Cells.Item(165, 4).Value = "Bathroom"
$ws.Cells.Item(165, 5).Value = "No Motion"
$ws.Cells.Item(165, 6).Value = "Inactive"
$ws.Cells.Item(166, 1).Value = "2026-01-28"
$ws.Cells.Item(166, 2).Value = "16:24:12"
$ws.Cells.Item(166, 3).Value = "16:00"
$ws.Cells.Item(166, 4).Value = "Bathroom"
$ws.Cells.Item(166, 5).Value = "No Motion"
$ws.Cells.Item(166, 6).Value = "Inactive"
$ws.Cells.Item(167, 1).Value = "2026-01-28"
$ws.Cells.Item(167, 2).Value = "16:24:17"
$ws.Cells.Item(167, 3).Value = "16:00"
$ws.Cells.Item(167, 4).Value = "Bathroom"
$ws.Cells.Item(167, 5).Value = "No Motion"
$ws.Cells.Item(167, 6).Value = "Inactive"
$ws.Cells.Item(168, 1).Value = "2026-01-28"
$ws.Cells.Item(168, 2).Value = "16:24:22"
$ws.Cells.Item(168, 3).Value = "16:00"
$ws.Cells.Item(168, 4).Value = "Bathroom"
$ws.Cells.Item(168, 5).Value = "No Motion"
$ws.Cells.Item(168, 6).Value = "Inactive"
$ws.Cells.Item(169, 1).Value = "2026-01-28"
$ws.Cells.Item(169, 2).Value = "16:24:27"
$ws.Cells.Item(169, 3).Value = "16:00"
$ws.Cells.Item(169, 4).Value = "Bathroom"
$ws.Cells.Item(169, 5).Value = "No Motion"
$ws.Cells.Item(169, 6).Value = "Inactive"
$ws.Cells.Item(170, 1).Value = "2026-01-28"
$ws.Cells.Item(170, 2).Value = "16:24:32"
$ws.Cells.Item(170, 3).Value = "16:00"
$ws.Cells.Item(170, 4).Value = "Bathroom"
$ws.Cells.Item(170, 5).Value = "No Motion"
$ws.Cells.Item(170, 6).Value = "Inactive"
$ws.Cells.Item(171, 1).Value = "2026-01-28"
$ws.Cells.Item(171, 2).Value = "16:24:37"
$ws.Cells.Item(171, 3).Value = "16:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "No Motion"
$ws.Cells.Item(171, 6).Value = "Inactive"
$ws.Cells.Item(172, 1).Value = "2026-01-28"
$ws.Cells.Item(172, 2).Value = "16:24:43"
$ws.Cells.Item(172, 3).Value = "16:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "No Motion"
$ws.Cells.Item(172, 6).Value = "Inactive"
$dateRng.Style = "Normal"

# ---- Humidity sheet: add rows 159-172 ----
$ws = $wb.Worksheets.Item("Humidity")
$dateRng = $ws.Range("A159:A172")
$dateRng.NumberFormat = "@"
$valRng = $ws.Range("E159:E172")
$valRng.NumberFormat = "@"
$ws.Cells.Item(159, 1).Value = "2026-01-28"
$ws.Cells.Item(159, 2).Value = "16:23:46"
$ws.Cells.Item(159, 3).Value = "16:00"
$ws.Cells.Item(159, 4).Value = "Bathroom"
$ws.Cells.Item(159, 5).Value = "87.4%"
$ws.Cells.Item(159, 6).Value = "Active"
$ws.Cells.Item(160, 1).Value = "2026-01-28"
$ws.Cells.Item(160, 2).Value = "16:23:49"
$ws.Cells.Item(160, 3).Value = "16:00"
$ws.Cells.Item(160, 4).Value = "Bathroom"
$ws.Cells.Item(160, 5).Value = "88.3%"
$ws.Cells.Item(160, 6).Value = "Active"
$ws.Cells.Item(161, 1).Value = "2026-01-28"
$ws.Cells.Item(161, 2).Value = "16:23:53"
$ws.Cells.Item(161, 3).Value = "16:00"
$ws.Cells.Item(161, 4).Value = "Bathroom"
$ws.Cells.Item(161, 5).Value = "87.4%"
$ws.Cells.Item(161, 6).Value = "Active"
$ws.Cells.Item(162, 1).Value = "2026-01-28"
$ws.Cells.Item(162, 2).Value = "16:23:57"
$ws.Cells.Item(162, 3).Value = "16:00"
$ws.Cells.Item(162, 4).Value = "Bathroom"
$ws.Cells.Item(162, 5).Value = "88.3%"
$ws.Cells.Item(162, 6).Value = "Active"
$ws.Cells.Item(163, 1).Value = "2026-01-28"
$ws.Cells.Item(163, 2).Value = "16:24:05"
$ws.Cells.Item(163, 3).Value = "16:00"
$ws.Cells.Item(163, 4).Value = "Bathroom"
$ws.Cells.Item(163, 5).Value = "88.3%"
$ws.Cells.Item(163, 6).Value = "Active"
$ws.Cells.Item(164, 1).Value = "2026-01-28"
$ws.Cells.Item(164, 2).Value = "16:24:09"
$ws.Cells.Item(164, 3).Value = "16:00"
$ws.Cells.Item(164, 4).Value = "Bathroom"
$ws.Cells.Item(164, 5).Value = "88.3%"
$ws.Cells.Item(164, 6).Value = "Active"
$ws.Cells.Item(165, 1).Value = "2026-01-28"
$ws.Cells.Item(165, 2).Value = "16:24:13"
$ws.Cells.Item(165, 3).Value = "16:00"
$ws.Cells.Item(165, 4).Value = "Bathroom"
$ws.Cells.Item(165, 5).Value = "87.4%"
$ws.Cells.Item(165, 6).Value = "Active"
$ws.Cells.Item(166, 1).Value = "2026-01-28"
$ws.Cells.Item(166, 2).Value = "16:24:18"
$ws.Cells.Item(166, 3).Value = "16:00"
$ws.Cells.Item(166, 4).Value = "Bathroom"
$ws.Cells.Item(166, 5).Value = "88.3%"
$ws.Cells.Item(166, 6).Value = "Active"
$ws.Cells.Item(167, 1).Value = "2026-01-28"
$ws.Cells.Item(167, 2).Value = "16:24:21"
$ws.Cells.Item(167, 3).Value = "16:00"
$ws.Cells.Item(167, 4).Value = "Bathroom"
$ws.Cells.Item(167, 5).Value = "88.3%"
$ws.Cells.Item(167, 6).Value = "Active"
$ws.Cells.Item(168, 1).Value = "2026-01-28"
$ws.Cells.Item(168, 2).Value = "16:24:26"
$ws.Cells.Item(168, 3).Value = "16:00"
$ws.Cells.Item(168, 4).Value = "Bathroom"
$ws.Cells.Item(168, 5).Value = "87.4%"
$ws.Cells.Item(168, 6).Value = "Active"
$ws.Cells.Item(169, 1).Value = "2026-01-28"
$ws.Cells.Item(169, 2).Value = "16:24:30"
$ws.Cells.Item(169, 3).Value = "16:00"
$ws.Cells.Item(169, 4).Value = "Bathroom"
$ws.Cells.Item(169, 5).Value = "88.3%"
$ws.Cells.Item(169, 6).Value = "Active"
$ws.Cells.Item(170, 1).Value = "2026-01-28"
$ws.Cells.Item(170, 2).Value = "16:24:34"
$ws.Cells.Item(170, 3).Value = "16:00"
$ws.Cells.Item(170, 4).Value = "Bathroom"
$ws.Cells.Item(170, 5).Value = "87.4%"
$ws.Cells.Item(170, 6).Value = "Active"
$ws.Cells.Item(171, 1).Value = "2026-01-28"
$ws.Cells.Item(171, 2).Value = "16:24:38"
$ws.Cells.Item(171, 3).Value = "16:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "88.3%"
$ws.Cells.Item(171, 6).Value = "Active"
$ws.Cells.Item(172, 1).Value = "2026-01-28"
$ws.Cells.Item(172, 2).Value = "16:24:42"
$ws.Cells.Item(172, 3).Value = "16:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "88.3%"
$ws.Cells.Item(172, 6).Value = "Active"
$dateRng.Style = "Normal"
$valRng.Style = "Normal"

# ---- Temperature sheet: add rows 159-172 ----
$ws = $wb.Worksheets.Item("Temperature")
$dateRng = $ws.Range("A159:A172")
$dateRng.NumberFormat = "@"
$ws.Cells.Item(159, 1).Value = "2026-01-28"
$ws.Cells.Item(159, 2).Value = "16:23:46"
$ws.Cells.Item(159, 3).Value = "16:00"
$ws.Cells.Item(159, 4).Value = "Bathroom"
$ws.Cells.Item(159, 5).Value = "22.8C"
$ws.Cells.Item(159, 6).Value = "Active"
$ws.Cells.Item(160, 1).Value = "2026-01-28"
$ws.Cells.Item(160, 2).Value = "16:23:50"
$ws.Cells.Item(160, 3).Value = "16:00"
$ws.Cells.Item(160, 4).Value = "Bathroom"
$ws.Cells.Item(160, 5).Value = "22.8C"
$ws.Cells.Item(160, 6).Value = "Active"
$ws.Cells.Item(161, 1).Value = "2026-01-28"
$ws.Cells.Item(161, 2).Value = "16:23:54"
$ws.Cells.Item(161, 3).Value = "16:00"
$ws.Cells.Item(161, 4).Value = "Bathroom"
$ws.Cells.Item(161, 5).Value = "22.8C"
$ws.Cells.Item(161, 6).Value = "Active"
$ws.Cells.Item(162, 1).Value = "2026-01-28"
$ws.Cells.Item(162, 2).Value = "16:23:58"
$ws.Cells.Item(162, 3).Value = "16:00"
$ws.Cells.Item(162, 4).Value = "Bathroom"
$ws.Cells.Item(162, 5).Value = "22.8C"
$ws.Cells.Item(162, 6).Value = "Active"
$ws.Cells.Item(163, 1).Value = "2026-01-28"
$ws.Cells.Item(163, 2).Value = "16:24:06"
$ws.Cells.Item(163, 3).Value = "16:00"
$ws.Cells.Item(163, 4).Value = "Bathroom"
$ws.Cells.Item(163, 5).Value = "22.8C"
$ws.Cells.Item(163, 6).Value = "Active"
$ws.Cells.Item(164, 1).Value = "2026-01-28"
$ws.Cells.Item(164, 2).Value = "16:24:10"
$ws.Cells.Item(164, 3).Value = "16:00"
$ws.Cells.Item(164, 4).Value = "Bathroom"
$ws.Cells.Item(164, 5).Value = "22.7C"
$ws.Cells.Item(164, 6).Value = "Active"
$ws.Cells.Item(165, 1).Value = "2026-01-28"
$ws.Cells.Item(165, 2).Value = "16:24:14"
$ws.Cells.Item(165, 3).Value = "16:00"
$ws.Cells.Item(165, 4).Value = "Bathroom"
$ws.Cells.Item(165, 5).Value = "22.8C"
$ws.Cells.Item(165, 6).Value = "Active"
$ws.Cells.Item(166, 1).Value = "2026-01-28"
$ws.Cells.Item(166, 2).Value = "16:24:18"
$ws.Cells.Item(166, 3).Value = "16:00"
$ws.Cells.Item(166, 4).Value = "Bathroom"
$ws.Cells.Item(166, 5).Value = "22.8C"
$ws.Cells.Item(166, 6).Value = "Active"
$ws.Cells.Item(167, 1).Value = "2026-01-28"
$ws.Cells.Item(167, 2).Value = "16:24:22"
$ws.Cells.Item(167, 3).Value = "16:00"
$ws.Cells.Item(167, 4).Value = "Bathroom"
$ws.Cells.Item(167, 5).Value = "22.8C"
$ws.Cells.Item(167, 6).Value = "Active"
$ws.Cells.Item(168, 1).Value = "2026-01-28"
$ws.Cells.Item(168, 2).Value = "16:24:26"
$ws.Cells.Item(168, 3).Value = "16:00"
$ws.Cells.Item(168, 4).Value = "Bathroom"
$ws.Cells.Item(168, 5).Value = "22.8C"
$ws.Cells.Item(168, 6).Value = "Active"
$ws.Cells.Item(169, 1).Value = "2026-01-28"
$ws.Cells.Item(169, 2).Value = "16:24:30"
$ws.Cells.Item(169, 3).Value = "16:00"
$ws.Cells.Item(169, 4).Value = "Bathroom"
$ws.Cells.Item(169, 5).Value = "22.7C"
$ws.Cells.Item(169, 6).Value = "Active"
$ws.Cells.Item(170, 1).Value = "2026-01-28"
$ws.Cells.Item(170, 2).Value = "16:24:34"
$ws.Cells.Item(170, 3).Value = "16:00"
$ws.Cells.Item(170, 4).Value = "Bathroom"
$ws.Cells.Item(170, 5).Value = "22.8C"
$ws.Cells.Item(170, 6).Value = "Active"
$ws.Cells.Item(171, 1).Value = "2026-01-28"
$ws.Cells.Item(171, 2).Value = "16:24:38"
$ws.Cells.Item(171, 3).Value = "16:00"
$ws.Cells.Item(171, 4).Value = "Bathroom"
$ws.Cells.Item(171, 5).Value = "22.7C"
$ws.Cells.Item(171, 6).Value = "Active"
$ws.Cells.Item(172, 1).Value = "2026-01-28"
$ws.Cells.Item(172, 2).Value = "16:24:42"
$ws.Cells.Item(172, 3).Value = "16:00"
$ws.Cells.Item(172, 4).Value = "Bathroom"
$ws.Cells.Item(172, 5).Value = "22.8C"
$ws.Cells.Item(172, 6).Value = "Active"
$dateRng.Style = "Normal"

$ws = $wb.Worksheets.Item("PIR")
$ws.Activate()
$ws.Range("A1").Select()